$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the s_vals sheet, matching the existing
# "sum" header's formatting (bold header style, border, centered).
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill the new column's data rows with 0 (default/placeholder save value).
$ws.Range("H2:H5").Value = 0
